$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 45 new rows (102-146) of test data following the existing repeating
# pattern (regcntr_id cycles 10002..10010, device_id increments, lang_code
# "eng", is_active TRUE, cr_by "superadmin", cr_dtimes "now()").
$aCycle = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)
$bStart = 3000121

for ($i = 0; $i -lt 45; $i++) {
    $r = 102 + $i
    $aVal = $aCycle[$i % $aCycle.Length]
    $bVal = $bStart + $i

    $ws.Cells.Item($r, 1).Value = $aVal
    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Update the view so the new rows are visible (matches the scroll/selection
# captured in the authored workbook).
[void]$ws.Range("A102:F146").Select()
$excel.ActiveWindow.ScrollRow = 128

# The authored workbook was printed/previewed at some point, which persists
# explicit page setup info (portrait orientation) to the sheet.
$ws.PageSetup.Orientation = 1
